# Generate Report for Handoff
#
# The localization pipeline picked up a new source file
# (41df4f83-74ed-45cb-97fd-e809bb4f37ad.md) superseding the prior one
# (bd18bc81-4910-4935-8288-b26f75fb1fb1.md). Refresh every sheet with the
# new file name / xliff hash, bump the "latest" timestamps, and clear the
# now-stale "Latest Target File" / "Latest Handback File" + its hyperlink
# on the per-locale sheets (handback hasn't happened yet for the new file).

$wb = $excel.ActiveWorkbook

$oldGuid = "bd18bc81-4910-4935-8288-b26f75fb1fb1"
$newGuid = "41df4f83-74ed-45cb-97fd-e809bb4f37ad"
$newHash = "6607bba356f414ec00275e023df1a5ac0679f8b0"

# Original external link (unchanged host/branch - only the md file name
# referenced by the commit message changes).
$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19997dfa071dfc176b7c0cc7398eaebd361656ae/e2e/$oldGuid.md"

# BGR-packed version of the custom Hyperlink font color (#6495ED) used
# throughout this workbook, for re-applying to re-created hyperlink cells.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-03 23:03:59"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $srcUrl, [System.Type]::Missing, [System.Type]::Missing, "e2e\$newGuid.md")
$wsOverview.Range("B2").Font.Color = $hyperlinkColor

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-09-03 23:03:54"
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("J2").Style = "Normal"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $srcUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$wsZh.Range("A2").Font.Color = $hyperlinkColor

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-09-03 23:03:59"
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("J2").Style = "Normal"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $srcUrl, [System.Type]::Missing, [System.Type]::Missing, "$newGuid.md")
$wsDe.Range("A2").Font.Color = $hyperlinkColor

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
